# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 10:05"

# --- Swap country rows whose shared-string order changed -------------
# Santa Lucia (row 201) <-> Belice (row 202): labels + D/H stats swap
$ws.Range("A201").Value = "Belice"
$ws.Range("D201").Value = 16
$ws.Range("H201").Value = 2

$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("D202").Value = 18
$ws.Range("H202").Value = 0

# Montserrat (row 210) <-> Seychelles (row 211): labels + D/H stats swap
$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# Bonaire, San Eustaquio y Saba (row 215) <-> San Bartolome (row 216)
$ws.Range("A215").Value = "San Bartolome"
$ws.Range("A216").Value = "Bonaire, San Eustaquio y Saba"

# --- Update COVID-19 stats for several countries ----------------------
# Rusia (row 6)
$ws.Range("B6").Value = 414878
$ws.Range("C6").Value = 9035
$ws.Range("D6").Value = 175877
$ws.Range("E6").Value = 234146
$ws.Range("G6").Value = 162
$ws.Range("H6").Value = 4855

# India (row 10)
$ws.Range("B10").Value = 190962
$ws.Range("C10").Value = 353
$ws.Range("D10").Value = 91866
$ws.Range("E10").Value = 93685
$ws.Range("G10").Value = 3
$ws.Range("H10").Value = 5411

# Singapur (row 29)
$ws.Range("B29").Value = 35292
$ws.Range("C29").Value = 408
$ws.Range("E29").Value = 13570

# Rumania (row 41)
$ws.Range("E41").Value = 4731
$ws.Range("G41").Value = 4
$ws.Range("H41").Value = 1270

# Afganistan (row 48)
$ws.Range("B48").Value = 15750
$ws.Range("C48").Value = 545
$ws.Range("D48").Value = 1428
$ws.Range("E48").Value = 14057
$ws.Range("G48").Value = 8
$ws.Range("H48").Value = 265

# Armenia (row 58)
$ws.Range("B58").Value = 9492
$ws.Range("C58").Value = 210
$ws.Range("E58").Value = 5951

# Chequia (row 60)
$ws.Range("B60").Value = 9273
$ws.Range("C60").Value = 5
$ws.Range("D60").Value = 6562
$ws.Range("E60").Value = 2391

# Estonia (row 95)
$ws.Range("B95").Value = 1870
$ws.Range("C95").Value = 1
$ws.Range("D95").Value = 1625

# Lituania (row 100)
$ws.Range("B100").Value = 1678
$ws.Range("C100").Value = 3
$ws.Range("E100").Value = 372

# Eslovaquia (row 103)
$ws.Range("B103").Value = 1522
$ws.Range("C103").Value = 1
$ws.Range("D103").Value = 1368
$ws.Range("E103").Value = 126
